# Use Case.docx - mini alteração no use case da eliminação de conta
# (and a handful of other run-merge only edits elsewhere in the doc)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "... O sistema  adiciona a transação ..." - merge the lone space
#    run with the following "adiciona a transação ..." run.
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute(
    "adiciona a transação ao histórico de transações do utilizador e envia o valor para o saldo do vendedor",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "adiciona a transação ao histórico de transações do utilizador e envia o valor para o saldo do vendedor",
    2)

# ---------------------------------------------------------------------
# 2) "3" + ".1" + ") O utilizador pesquisa..." -> single run "3.1) O utilizador pesquisa..."
#    "3.2" + ") O sistema verifica..."        -> single run "3.2) O sistema verifica..."
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute(
    ".1) O utilizador pesquisa a transação através de palavras-chave ou o ID da transação",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ".1) O utilizador pesquisa a transação através de palavras-chave ou o ID da transação",
    2)

$r = $d.Content
$null = $r.Find.Execute(
    ") O sistema verifica a pesquisa e exibe os resultados",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ") O sistema verifica a pesquisa e exibe os resultados",
    2)

# ---------------------------------------------------------------------
# 3) "Fluxo Alternativo (" + "2" + ")" -> single run "Fluxo Alternativo (2)"
#    "[" + "O utilizador não encontra..." + "] (passo " + "3.1" + ")"
#        -> single run "[O utilizador não encontra o que procura na pesquisa] (passo 3.1)"
# ---------------------------------------------------------------------
$anchor = $d.Content
$null = $anchor.Find.Execute(
    "3.2) O sistema verifica a pesquisa e exibe os resultados",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scope = $d.Range($anchor.End, $anchor.End + 200)
$null = $scope.Find.Execute("2)", $false, $false, $false, $false, $false, $true, 1, $false, "2)", 2)

$r = $d.Content
$null = $r.Find.Execute(
    "[O utilizador não encontra o que procura na pesquisa] (passo 3.1)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "[O utilizador não encontra o que procura na pesquisa] (passo 3.1)",
    2)

# ---------------------------------------------------------------------
# 4) "Fluxo " + "Alternativo" + " (1)" -> single run "Fluxo Alternativo (1)"
#    "[" + " " + "Os dados inseridos..." + " ou password antiga errada" + "] (passo " + "5" + ")"
#        -> single run "[ Os dados inseridos pelo utilizador eram duplicados ou password antiga errada] (passo 5)"
# ---------------------------------------------------------------------
$anchor = $d.Content
$null = $anchor.Find.Execute(
    "7) O utilizador conseguiu alterar os seus dados e vê-los no seu perfil",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scope = $d.Range($anchor.End, $anchor.End + 200)
$null = $scope.Find.Execute(
    "Fluxo Alternativo (1)", $false, $false, $false, $false, $false, $true, 1, $false,
    "Fluxo Alternativo (1)", 2)

$r = $d.Content
$target = "[ Os dados inseridos pelo utilizador eram duplicados ou password antiga errada] (passo 5)"
$null = $r.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, $target, 2)

# ---------------------------------------------------------------------
# 5) "Pedido de Eliminação da Conta" -> "Pedir Eliminação da Conta"
#    (runs: "Pedi" / "r" / " Eliminação da Conta")
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Pedido de Eliminação da Conta", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sub = $d.Range($r.Start + 4, $r.Start + 9)
$sub.Text = "r"
$sub.Font.Bold = 1
$sub.Font.Bold = 0

# ---------------------------------------------------------------------
# 6) "Processo de pedir para eliminar a tua conta" -> "...a sua conta"
#    (runs: "Processo de pedir para eliminar a " / "s" / "ua conta")
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Processo de pedir para eliminar a tua conta", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sub = $d.Range($r.Start + 34, $r.Start + 35)
$sub.Text = "s"
$sub.Font.Bold = 1
$sub.Font.Bold = 0
